# Rename the embedded logo pictures in the document's headers/footers.
#
#  - Header (first page): BTec_Logo-Orange picture  image2.jpg -> image1.jpg
#  - Footer (default):    Pearson logo picture       image1.png -> image2.png
#  - Footer (first page): Pearson logo picture       image1.png -> image2.png
#
# Each picture is an InlineShape; we rename it by setting InlineShape.Name,
# which is the COM-interop equivalent of renaming the picture in the Word UI
# (Alt Text / picture name). We re-fetch each shape through the paragraph
# that directly contains it (rather than through the whole header/footer
# Range) so the object handle addresses the right block before we mutate it.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-LogoInHeaderFooter($story, $newName) {
    $rng = $story.Range
    $paraCount = $rng.Paragraphs.Count
    $para = $rng.Paragraphs($paraCount)
    $shape = $para.Range.InlineShapes(1)
    $shape.Name = $newName
}

# --- Header: first-page header holds the BTec logo ---
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers($i)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        Rename-LogoInHeaderFooter $hdr "image1.jpg"
    }
}

# --- Footers: both default and first-page footers hold the Pearson logo ---
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers($i)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        Rename-LogoInHeaderFooter $ftr "image2.png"
    }
}
